$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new worksheet "A06 vie saint gregoire". It becomes the
#    6th tab, pushing every sheet from the old #6 onward one slot to
#    the right (handled by the renumbering pass below).
# ------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(6)
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "A06 vie saint gregoire"

$header = @("line_n", "prev_line", "line", "next_line")
for ($c = 1; $c -le 4; $c++) {
    $newSheet.Cells.Item(1, $c).Value = $header[$c - 1]
}
$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$data = @(
    @(37, "Quant le damoiseau a celle nouvelle ouÿe,", "Grant dueil eult en son cueur, puis dit: “Virge marie,", "De si vilain pechié ne nous souvient il mie.”"),
    @(39, "De si vilain pechié ne nous souvient il mie.”", "Puis appella sa seur et lui dist: “Doulce amie,", "Ces gens de dehors dient, c'est chose prouvee,"),
    @(51, "Du damoiseau dirons, qui fut courtois et sage:", "Il a dit a sa seur: “Vecy trop grant servage,", "Car tous ceulx sont dampnez en mon entencïon"),
    @(68, "Mais quant il le regarde, moult oult le cueur dolent.", "Puis lui dist: “Belle seur, que pourrons devenir?", "Pour parolle de gens, ne laissons a mourir."),
    @(71, "Les deux cueurs de noz ventres deussent parmi partir.”", "Et elle respondi: “Dieu fera son plaisir.", "Faittes faire, beau frere, ung escrin sans delay,"),
    @(76, "Comment vous l'engendrastes, et comment le portay.”", "Son frere respondi: “Je feray vo talent.”", "Ung panier a fait faire, bien tost et promptement."),
    @(152, "Et cil dit qu'il l'ara a mouller et a per,", "Et tousjours luy respont: “Ja ne m'espousera.”", "Et le conte lui dit tousiours que si fera."),
    @(160, "Ung jour par ung matin, s'en vint devant l'abbé,", "Et lui dist sagement: “Sire, qui est mon pere?", "Et me dittes aussi nouvelles de ma mere.”"),
    @(199, "Chacun si lui demande de quel terre il est né.", "Dist Gregoire: “Je suis d'un estrange contre,", "Et qui a ceste terre si laidement gastee.”"),
    @(208, "“Damoiseau,” dist la dame, “de quel terre estes né?”", "Dist Gregoire: “Je suis d'un estrange contré.”", "Vray Dieu, elle l'avoit en son ventre porté,"),
    @(216, "Lors respondy Gregoire sans soy point delayer,", "Et a dit a la dame: “Je ferai vo tallent.", "Faittes mander au conte sans nul delayement,"),
    @(249, "Gregoire a abbatu le conte du cheval;", "La dame qui le vit dist: “C'est coup de vassal!", "Je vous pry de bon cueur, doulz Pere espirital,"),
    @(270, "Saches certainement que moult fut esperdu.", "Il a dit a Gregoire: “Je me tieng pour vaincu.", "Pour dieu, ne m'ochis mie, qui en croix fut pendu.”"),
    @(295, "Et de sa gentil dame qui tant ot cler le vis.", "Doulcement luy a dit: “Vous estes mes amis", "Gentilz homs,” dit la dame, “s'espouser me voulez,"),
    @(321, "Sa femme et sa mere l'aimoit de grant amour.", "Un soir apres souper lui dist: “Mon chier seignour,", "Je vous perchoy trop fort pallir vostre coulour."),
    @(339, "Or commenche a gregoire sa paine et son tourment.", "Si tost come le sceüt, il a dit: “Dieu, vray pere,", "Or suis je engendrez en la seur et du frere;"),
    @(344, "Sa mere vint a luy, comme femme desvee,", "En disant: “Mon enfant, quel povre destinee!", "En mes flans vous portay, et m'avez espousee.”"),
    @(355, "Pour leur seigneur Gregoire, qui s'en vouloit aler.", "Sa mere lui dist: “Filz, avecques vous yray.”", "“Non ferez, dame, par la foy que vous doy."),
    @(373, "Et sa mere demeure, qui forment se demente,", "Et dist: “Las, que feray, la chetive dolente,", "Quat mon filz ma laisse mort vien et si mavance?"),
    @(397, "Son estat lui a dit Gregoire maintenant.", "Et lui dist: “Je vouldroye avoir herbergement", "Auquel peusse servir mon Dieu benignement."),
    @(409, "Qu'il auoit une roche dedens ysle de mer,", "Il lui dist doulcement: “Se m'y voulez mener,", "Mon cheval et ma robe vous vouldroye doner,"),
    @(412, "Et trestout mon argent aussi vous donneray.”", "L'oste luy respondi: “Et je vous y menray.", "Dittes moy vostre nom, s'il vous plaist, sans delay.”"),
    @(430, "Si tost com il la tint, en la mer la jetta,", "Puis a dit a son hoste: “Je demouray decha.", "Jamais de cy n'yray, si sera retrouvee.”"),
    @(495, "Dieu veult que pape soyes a Rome sans delay.”", "Saint Gregoire respond: “Par la foy qu'a Dieu doy,", "Jamais de ceste roche mie ne partiray,"),
    @(500, "“Beau sire veez la, n'en soyez en doubtance.”", "Quant il la vit si dist: “Dieu, qu'avez grant puissace!", "Or voy je bien qu'ay fait de mes maulx penitace."),
    @(549, "“Pere, plus grant pecheresse ne fut onc mais ouÿe.”", "Saint Gregoire l'apelle, et lui dist: “Doulce amye,", "Contez moy vostre estat maintenant, sans delay.”"),
    @(551, "Contez moy vostre estat maintenant, sans delay.”", "Elle dist: “De mon frere oulz ung filz, par ma foy.", "Encores puis apres, mon enffant espousay."),
    @(558, "L'absolut et puis print la besongne a compter,", "Et dist: “Vous estes celle qui neuf mois me porta.”", "Quant elle l'entendi, grans joye en demena.")
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Renumber the "A##" prefix of every sheet from the old #6 onward
#    by +1, to make room for the newly inserted sheet above. Sheets
#    are renamed from the last tab backwards so no two sheets ever
#    collide on the same name mid-script.
# ------------------------------------------------------------------
$wb.Worksheets.Item(23).Name = "A26 quatre sereurs"
$wb.Worksheets.Item(22).Name = "A24 roy avoit amie"
$wb.Worksheets.Item(21).Name = "A23 jeu des dez"
$wb.Worksheets.Item(20).Name = "A22 mauvais riche homme"
$wb.Worksheets.Item(19).Name = "A21 vieillards tués"
$wb.Worksheets.Item(18).Name = "A20 elegy troyes"
$wb.Worksheets.Item(17).Name = "A19 richart sans peour"
$wb.Worksheets.Item(16).Name = "A18 robert deable"
$wb.Worksheets.Item(15).Name = "A17 guillaume angleterre"
$wb.Worksheets.Item(14).Name = "A14 vie seint thibault"
$wb.Worksheets.Item(13).Name = "A13 miracle saint servais"
$wb.Worksheets.Item(12).Name = "A12 vie saint sebastien"
$wb.Worksheets.Item(11).Name = "A11 poines enfer"
$wb.Worksheets.Item(10).Name = "A10 vie saint leu"
$wb.Worksheets.Item(9).Name = "A09 vie glorieux confesseur"
$wb.Worksheets.Item(8).Name = "A08 vie saint jean paulus"
$wb.Worksheets.Item(7).Name = "A07 saint jean evangeliste"
